$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1175
$ws.Range("J127").Value = 1175
$ws.Range("L127").Value = 3525
$ws.Range("N127").Value = -13445
$ws.Range("H135").Value = 853.3333
$ws.Range("I135").Value = 576.6667
$ws.Range("J135").Value = 1545
$ws.Range("K135").Value = 5190.0003
$ws.Range("L135").Value = 13905
$ws.Range("M135").Value = -2655.0003
$ws.Range("N135").Value = -18975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 38000
$ws.Range("J22").Value = 38000
$ws.Range("L22").Value = 38000
$ws.Range("N22").Value = -38598
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 18005.334
$ws.Range("J27").Value = 18005.334
$ws.Range("L27").Value = 18005.334
$ws.Range("N27").Value = -18373.334
$ws.Range("H43").Value = 25301.6
$ws.Range("J43").Value = 25301.6
$ws.Range("L43").Value = 25301.6
$ws.Range("N43").Value = -25927.6
$ws.Range("H109").Value = 18995
$ws.Range("J109").Value = 18995
$ws.Range("L109").Value = 18995
$ws.Range("N109").Value = -21769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -2340
$ws.Range("H23").Value = 23511.945
$ws.Range("I23").Value = 9061.799999999999
$ws.Range("J23").Value = 29069.691
$ws.Range("K23").Value = 9061.799999999999
$ws.Range("L23").Value = 29069.691
$ws.Range("M23").Value = -8821.799999999999
$ws.Range("N23").Value = -29549.691
$ws.Range("H27").Value = 23511.945
$ws.Range("I27").Value = 9061.799999999999
$ws.Range("J27").Value = 29069.691
$ws.Range("K27").Value = 9061.799999999999
$ws.Range("L27").Value = 29069.691
$ws.Range("M27").Value = -8869.799999999999
$ws.Range("N27").Value = -29453.691
$ws.Range("H107").Value = 602.4211
$ws.Range("I107").Value = 246.72728
$ws.Range("J107").Value = 1091.5
$ws.Range("K107").Value = 246.72728
$ws.Range("L107").Value = 1091.5
$ws.Range("M107").Value = 1673.27272
$ws.Range("N107").Value = -4931.5
$ws.Range("H122").Value = 30304082
$ws.Range("I122").Value = 47619824
$ws.Range("J122").Value = 1537.5
$ws.Range("K122").Value = 142859472
$ws.Range("L122").Value = 4612.5
$ws.Range("M122").Value = -142857022
$ws.Range("N122").Value = -9512.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 218.16667
$ws.Range("I17").Value = 136.33333
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 408.99999
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = -239.99999
$ws.Range("N17").Value = -1238
$ws.Range("H34").Value = 2131.6667
$ws.Range("J34").Value = 2458
$ws.Range("L34").Value = 7374
$ws.Range("N34").Value = -7542
$ws.Range("H39").Value = 977.3333
$ws.Range("J39").Value = 977.3333
$ws.Range("L39").Value = 2931.9999
$ws.Range("N39").Value = -3519.9999
$ws.Range("H55").Value = 626
$ws.Range("I55").Value = 304
$ws.Range("J55").Value = 733.3333
$ws.Range("K55").Value = 912
$ws.Range("L55").Value = 2199.9999
$ws.Range("M55").Value = -735
$ws.Range("N55").Value = -2553.9999
$ws.Range("H107").Value = 41667390
$ws.Range("J107").Value = 790.9091
$ws.Range("L107").Value = 2372.7273
$ws.Range("N107").Value = -6212.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2311283
$ws.Range("I24").Value = 2887857.8
$ws.Range("J24").Value = 4984.5
$ws.Range("K24").Value = 2887857.8
$ws.Range("L24").Value = 4984.5
$ws.Range("M24").Value = -2887684.8
$ws.Range("N24").Value = -5330.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1913.7037
$ws.Range("I7").Value = 1799.1666
$ws.Range("K7").Value = 1799.1666
$ws.Range("M7").Value = -1687.1666
$ws.Range("H14").Value = 70003.75
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 70003.75
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 70003.75
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -70347.75
$ws.Range("H22").Value = 1264.4
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1480.5
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1480.5
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -2070.5
$ws.Range("H27").Value = 1264.4
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1480.5
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1480.5
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1694.5
$ws.Range("H31").Value = 753.75
$ws.Range("I31").Value = 838.3333
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 838.3333
$ws.Range("L31").Value = 500
$ws.Range("M31").Value = -590.3333
$ws.Range("N31").Value = -996
$ws.Range("H40").Value = 27030828
$ws.Range("I40").Value = 40002804
$ws.Range("J40").Value = 5873.75
$ws.Range("K40").Value = 40002804
$ws.Range("L40").Value = 5873.75
$ws.Range("M40").Value = -40002668
$ws.Range("N40").Value = -6145.75
$ws.Range("H126").Value = 1913.7037
$ws.Range("I126").Value = 1799.1666
$ws.Range("K126").Value = 5397.4998
$ws.Range("M126").Value = -2927.4998
$ws.Range("H136").Value = 20002198
$ws.Range("I136").Value = 33334372
$ws.Range("J136").Value = 3939
$ws.Range("K136").Value = 100003116
$ws.Range("L136").Value = 11817
$ws.Range("M136").Value = -100000566
$ws.Range("N136").Value = -16917

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 23316.75
$ws.Range("J109").Value = 23316.75
$ws.Range("L109").Value = 23316.75
$ws.Range("N109").Value = -26090.75
$ws.Range("H122").Value = 3204.9375
$ws.Range("I122").Value = 2532.348
$ws.Range("J122").Value = 4923.778
$ws.Range("K122").Value = 7597.044
$ws.Range("L122").Value = 14771.334
$ws.Range("M122").Value = -5147.044
$ws.Range("N122").Value = -19671.334
$ws.Range("H132").Value = 4112.067
$ws.Range("I132").Value = 3468.3
$ws.Range("J132").Value = 5399.6
$ws.Range("K132").Value = 10404.9
$ws.Range("L132").Value = 16198.8
$ws.Range("M132").Value = -7874.900000000001
$ws.Range("N132").Value = -21258.8
